# "Rule factoring" tracker: update a handful of rule-status cells and the
# wording of the "numeric tool version" rule note, per the commit
# "Update version-related comment in rules spreadsheet."
#
# Net effect on the status column (F) and the Notes column (H):
#   F5  (FileUrisMustConformToRfc8089)        TODO          -> IN PROGRESS
#   F6  (FileUrisMustNotIncludeDotDotSegments) DONE         -> IN PROGRESS
#   F22 (RegionStartPropertyMustBePresent)     DONE         -> TODO
#   H35 (MessagePropertiesMustBeConsistent)    note removed ("Talk to MF...")
#   F37 (IncludeDynamicContent)                DONE         -> IN PROGRESS
#   F38 (EnquoteDynamicContent)                DONE         -> IN PROGRESS
#   F49 (EliminateLocationOnlyArtifacts)       IN PROGRESS: EN -> TODO
#   H54 (UseNumericToolVersions)               comment text edited (see below)
#   F62 (UseConventionalRuleIds)               IN PROGRESS: EN -> TODO
#   F63 (UseConventionalUriBaseIdNames)        IN PROGRESS: EN -> TODO

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-Status($cellAddr, $templateAddr, $text) {
    $ws.Range($templateAddr).Copy() | Out-Null
    $ws.Range($cellAddr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($cellAddr).Value = $text
}

# Existing cells used purely as formatting templates (their own
# content/style is untouched by this change):
#   F8   -> the red "TODO" fill/style
#   G20  -> the gold "IN PROGRESS" fill/style

Set-Status "F5"  "G20" "IN PROGRESS"
Set-Status "F6"  "G20" "IN PROGRESS"
Set-Status "F22" "F8"  "TODO"
Set-Status "F37" "G20" "IN PROGRESS"
Set-Status "F38" "G20" "IN PROGRESS"
Set-Status "F49" "F8"  "TODO"
Set-Status "F62" "F8"  "TODO"
Set-Status "F63" "F8"  "TODO"

$excel.CutCopyMode = $false

# Drop the stale "talk to MF" note entirely (cell + formatting).
$ws.Range("H35").Clear()

# Reword the semanticVersion-related rule note.
$ws.Range("H54").Value = "Require an integer, followed by anything at all. We _don't_ validate semanticVersion, which should be a separate check in the same rule, and which in fact could have been in the schema!"

# Sheet view: reset the zoomed-in/scrolled state and move the selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.Zoom = 100
$ws.Range("H55").Select() | Out-Null
